# Dashboard de KPIs: Análisis de Lead Time y Gastos por demoras
#
# 1) Reorder sheets so "Dashboard" moves from first to last position:
#    importaciones, Packing_List, Auditoria_Facturacion, Dashboard
# 2) Add a new calculated column "Gastos de demora (USD)" to the
#    Tbl_Importaciones table (column J), with formula:
#    =IF([Días de despacho]>4,([Días de despacho]-4)*50,0)

$wb = $excel.ActiveWorkbook

# --- 1. Move Dashboard sheet to the end ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$dashboard.Move($null, $lastSheet)

# --- 2. Add new column to Tbl_Importaciones table ---
$ws = $wb.Worksheets.Item("importaciones")
$lo = $ws.ListObjects.Item("Tbl_Importaciones")
$newCol = $lo.ListColumns.Add()

# Set header text for the new column
$lo.HeaderRowRange.Cells.Item(1, $newCol.Index).Value2 = "Gastos de demora (USD)"

# Fill in the calculated column formula for each data row
$formula = "=IF(Tbl_Importaciones[[#This Row],[Días de despacho]]>4,(Tbl_Importaciones[[#This Row],[Días de despacho]]-4)*50,0)"
foreach ($r in $lo.DataBodyRange.Rows) {
    $r.Cells.Item(1, $newCol.Index).Formula = $formula
}

# Match the selection left on the "importaciones" sheet after adding the column
$ws.Activate() | Out-Null
$ws.Range("J3").Select() | Out-Null

$wb.Save()
